$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hourly-trigger time slot for this block moved from 18:20-18:29 to 22:00-22:09.
$ws.Range("B8").Value = "22:00 - 22:04"
$ws.Range("B9").Value = "22:05 - 22:09"

# Reflect the author's last on-screen selection/scroll position.
$ws.Range("A12").Select()
